# Auto-generated edit script applying numeric updates to leve-profit sheets.
# Source: diff of Sheets/Alexander_Profits.xlsx (scheduled price-refresh run).
$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 796.58826
$ws.Range("I33").Value = 1050.1818
$ws.Range("J33").Value = 331.66666
$ws.Range("K33").Value = 1050.1818
$ws.Range("L33").Value = 331.66666
$ws.Range("M33").Value = -821.1818000000001
$ws.Range("N33").Value = -789.66666
$ws.Range("H124").Value = 28987.5
$ws.Range("J124").Value = 28987.5
$ws.Range("L124").Value = 28987.5
$ws.Range("N124").Value = -38807.5
$ws.Range("H130").Value = 36830
$ws.Range("J130").Value = 36830
$ws.Range("L130").Value = 36830
$ws.Range("N130").Value = -46870
$ws.Range("H132").Value = 2402.0513
$ws.Range("I132").Value = 1804.3
$ws.Range("J132").Value = 4394.5557
$ws.Range("K132").Value = 5412.9
$ws.Range("L132").Value = 13183.6671
$ws.Range("M132").Value = -2882.9
$ws.Range("N132").Value = -18243.6671
$ws.Range("H137").Value = 2716.122
$ws.Range("I137").Value = 946.4761999999999
$ws.Range("J137").Value = 4574.25
$ws.Range("K137").Value = 2839.4286
$ws.Range("L137").Value = 13722.75
$ws.Range("M137").Value = -289.4285999999997
$ws.Range("N137").Value = -18822.75
$ws.Range("H138").Value = 2266.48
$ws.Range("I138").Value = 1055.725
$ws.Range("J138").Value = 3073.65
$ws.Range("K138").Value = 3167.175
$ws.Range("L138").Value = 9220.950000000001
$ws.Range("M138").Value = 1972.825
$ws.Range("N138").Value = -19500.95

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37002.832
$ws.Range("I32").Value = 13582.6045
$ws.Range("J32").Value = 137709.8
$ws.Range("K32").Value = 13582.6045
$ws.Range("L32").Value = 137709.8
$ws.Range("M32").Value = -13295.6045
$ws.Range("N32").Value = -138283.8
$ws.Range("H74").Value = 21186.809
$ws.Range("I74").Value = 1272.3903
$ws.Range("K74").Value = 1272.3903
$ws.Range("M74").Value = -398.3903
$ws.Range("H77").Value = 21186.809
$ws.Range("I77").Value = 1272.3903
$ws.Range("K77").Value = 6361.9515
$ws.Range("M77").Value = -1993.9515
$ws.Range("H133").Value = 39324.285
$ws.Range("J133").Value = 39324.285
$ws.Range("L133").Value = 39324.285
$ws.Range("N133").Value = -44384.285
$ws.Range("H135").Value = 47571.168
$ws.Range("J135").Value = 47571.168
$ws.Range("L135").Value = 47571.168
$ws.Range("N135").Value = -57711.168

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents() | Out-Null
$ws.Range("H124").Value = 42480
$ws.Range("J124").Value = 42480
$ws.Range("L124").Value = 42480
$ws.Range("N124").Value = -52300

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 30000
$ws.Range("J26").Value = 30000
$ws.Range("L26").Value = 30000
$ws.Range("N26").Value = -30574
$ws.Range("H50").Value = 9239.166999999999
$ws.Range("J50").Value = 9239.166999999999
$ws.Range("L50").Value = 9239.166999999999
$ws.Range("N50").Value = -10489.167
$ws.Range("H51").Value = 9234.666999999999
$ws.Range("J51").Value = 9234.666999999999
$ws.Range("L51").Value = 9234.666999999999
$ws.Range("N51").Value = -10706.667
$ws.Range("H61").Value = 9234.666999999999
$ws.Range("J61").Value = 9234.666999999999
$ws.Range("L61").Value = 9234.666999999999
$ws.Range("N61").Value = -9930.666999999999
$ws.Range("H97").Value = 30099
$ws.Range("J97").Value = 30099
$ws.Range("L97").Value = 30099
$ws.Range("N97").Value = -32081
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents() | Out-Null
$ws.Range("H130").Value = 56285
$ws.Range("J130").Value = 56285
$ws.Range("L130").Value = 56285
$ws.Range("N130").Value = -66325
$ws.Range("H134").Value = 4894.5483
$ws.Range("I134").Value = 5276.107
$ws.Range("K134").Value = 15828.321
$ws.Range("M134").Value = -13293.321

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1988.8889
$ws.Range("J92").Value = 1988.8889
$ws.Range("L92").Value = 5966.6667
$ws.Range("N92").Value = -8462.6667
$ws.Range("H107").Value = 1038.2
$ws.Range("I107").Value = 737.4286
$ws.Range("K107").Value = 2212.2858
$ws.Range("M107").Value = -292.2857999999997

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2837.5557
$ws.Range("J43").Value = 9619
$ws.Range("L43").Value = 9619
$ws.Range("N43").Value = -9921
$ws.Range("H93").Value = 12395.77
$ws.Range("J93").Value = 12395.77
$ws.Range("L93").Value = 12395.77
$ws.Range("N93").Value = -16139.77
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920
$ws.Range("H128").Value = 45780
$ws.Range("J128").Value = 45780
$ws.Range("L128").Value = 45780
$ws.Range("N128").Value = -55740
$ws.Range("H130").Value = 47686.668
$ws.Range("J130").Value = 47686.668
$ws.Range("L130").Value = 47686.668
$ws.Range("N130").Value = -57726.668
$ws.Range("H133").Value = 35681.25
$ws.Range("J133").Value = 35681.25
$ws.Range("L133").Value = 35681.25
$ws.Range("N133").Value = -45801.25
$ws.Range("H135").Value = 66772.22
$ws.Range("J135").Value = 66772.22
$ws.Range("L135").Value = 66772.22
$ws.Range("N135").Value = -76912.22

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 489
$ws.Range("I55").Value = 148.70589
$ws.Range("J55").Value = 971.0833
$ws.Range("K55").Value = 148.70589
$ws.Range("L55").Value = 971.0833
$ws.Range("M55").Value = 24.29410999999999
$ws.Range("N55").Value = -1317.0833
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents() | Out-Null
$ws.Range("H134").Value = 32770
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 32770
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 32770
$ws.Range("M134").ClearContents() | Out-Null
$ws.Range("N134").Value = -42910

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 36710.465
$ws.Range("I123").Value = 25000
$ws.Range("J123").Value = 44517.445
$ws.Range("K123").Value = 25000
$ws.Range("L123").Value = 44517.445
$ws.Range("M123").Value = -20100
$ws.Range("N123").Value = -54317.445
$ws.Range("H135").Value = 51184.645
$ws.Range("J135").Value = 51184.645
$ws.Range("L135").Value = 51184.645
$ws.Range("N135").Value = -61324.645
